$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the per-observation values between row 2 and row 3.
# Column A (Id), Q (Ost), R (Nord) are numeric.
# Column I (Antal) and AC (Publik kommentar) are stored as text.

$a2 = $ws.Range("A2").Value2
$a3 = $ws.Range("A3").Value2
$ws.Range("A2").Value2 = $a3
$ws.Range("A3").Value2 = $a2

# Column I is stored as text in the source data (e.g. "100", "15"),
# even though the text looks numeric. Use a leading apostrophe so Excel
# keeps these as text cells instead of auto-converting to numbers, then
# reset the style so we don't leave a visible "quote prefix" formatting
# behind on the cell.
$i2 = $ws.Range("I2").Value2
$i3 = $ws.Range("I3").Value2
$ws.Range("I2").Value2 = "'" + $i3
$ws.Range("I3").Value2 = "'" + $i2
$ws.Range("I2").Style = "Normal"
$ws.Range("I3").Style = "Normal"

$q2 = $ws.Range("Q2").Value2
$q3 = $ws.Range("Q3").Value2
$ws.Range("Q2").Value2 = $q3
$ws.Range("Q3").Value2 = $q2

$r2 = $ws.Range("R2").Value2
$r3 = $ws.Range("R3").Value2
$ws.Range("R2").Value2 = $r3
$ws.Range("R3").Value2 = $r2

$ac2 = $ws.Range("AC2").Value2
$ac3 = $ws.Range("AC3").Value2
$ws.Range("AC2").Value2 = $ac3
$ws.Range("AC3").Value2 = $ac2
